$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.889.51'
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').Value = '2.362.62'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('E5').Value = '  -1.47%  '
$ws.Range('D6').Value = '''239.86'
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').Value = '''74.30'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +0.24%  '
$ws.Range('E10').Value = '  +1.25%  '
$ws.Range('D11').Value = '''60.08'
$ws.Range('E11').Value = '  +4.86%  '
$ws.Range('D12').Value = '''37.14'
$ws.Range('E12').Value = '  +14.24%  '
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').Value = '''16.34'
$ws.Range('E15').Value = '  -1.67%  '
$ws.Range('D16').Value = '''0.927'
$ws.Range('E16').Value = '  +2.80%  '
$ws.Range('D17').Value = '2.371.32'
$ws.Range('E17').Value = '  +0.45%  '
$ws.Range('D18').Value = '43.846.20'
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('E19').Value = '  +1.24%  '
$ws.Range('E20').Value = '  -2.43%  '
$ws.Range('D21').Value = '''77.47'
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('D22').Value = '''253.77'
$ws.Range('E22').Value = '  -2.14%  '
$ws.Range('B23').Value = 'WEMIXToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D23').Value = '''3.80'
$ws.Range('E23').Value = '  +3.91%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = '''0.999'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  -4.22%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = '''10.55'
$ws.Range('E27').Value = '  -1.87%  '
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('E29').Value = '  -1.60%  '
$ws.Range('D30').Value = '''175.38'
$ws.Range('E30').Value = '  -0.58%  '
$ws.Range('E31').Value = '  +0.38%  '
$ws.Range('E32').Value = '  -1.51%  '
$ws.Range('D33').Value = '''0.0760'
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('D34').Value = '''5.43'
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('D35').Value = '''5.09'
$ws.Range('E35').Value = '  -3.05%  '
$ws.Range('D36').Value = '''3.79'
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('D37').Value = '''6.62'
$ws.Range('E37').Value = '  +4.67%  '
$ws.Range('E38').Value = '  +1.85%  '
$ws.Range('D39').Value = '''0.0281'
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('D40').Value = '''5.65'
$ws.Range('E40').Value = '  +20.25%  '
$ws.Range('D41').Value = '''20.48'
$ws.Range('E41').Value = '  +7.75%  '
$ws.Range('D42').Value = '''65.46'
$ws.Range('E42').Value = '  +12.07%  '
$ws.Range('E43').Value = '  -2.94%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''9.07'
$ws.Range('E44').Value = '  +0.91%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = '''0.203'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('D46').Value = '''2.55'
$ws.Range('E46').Value = '  +1.53%  '
$ws.Range('B47').Value = 'BinanceUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D47').Value = '''1.00'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').Value = '''1.24'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('D50').Value = '''98.40'
$ws.Range('E50').Value = '  -2.07%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = '''4.39'
$ws.Range('E51').Value = '  +14.90%  '

Write-Output "done"